{"js": "// Add two new list items after \"Building About Us - Part 02\":\n//   \"Building Testimonials - Part 01\"\n//   \"Building Testimonials - Part 02\"  (typed as \"...Part 0\" + \"2\")\n// Both paragraphs must inherit the same list/paragraph formatting\n// (ListParagraph style, numId=2, shading, spacing, outline level, fonts)\n// as the existing \"Building About Us\" bullet items right above them.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph (\"Building About Us - Part 02\") -- it is the\n// last paragraph in the document body.\nconst items = paragraphs.items;\nconst anchor = items[items.length - 1];\nanchor.load(\"text\");\nawait context.sync();\n\nif (anchor.text.trim() !== \"Building About Us - Part 02\") {\n  throw new Error(\n    \"Unexpected document structure: last paragraph is not \" +\n      '\"Building About Us - Part 02\" (found \"' + anchor.text + '\")'\n  );\n}\n\n// Insert the first new bullet right after the anchor paragraph. Office.js\n// copies the anchor's paragraph formatting (style/numbering/shading/\n// spacing/outline level/run fonts) onto the newly created paragraph.\nconst testimonials01 = anchor.insertParagraph(\n  \"Building Testimonials - Part 01\",\n  Word.InsertLocation.after\n);\n\n// Insert the second new bullet right after the first one, same formatting.\nconst testimonials02 = testimonials01.insertParagraph(\n  \"Building Testimonials - Part 0\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// The final \"2\" of \"Part 02\" was typed as its own run in the source\n// document, so append it as a separate insert on the paragraph's end range.\nconst endRange = testimonials02.getRange(Word.RangeLocation.end);\nendRange.insertText(\"2\", Word.InsertLocation.end);\nawait context.sync();\n\ntestimonials02.load(\"text\");\nawait context.sync();\n\nif (testimonials02.text !== \"Building Testimonials - Part 02\") {\n  throw new Error(\n    'Failed to build \"Building Testimonials - Part 02\" paragraph (found \"' +\n      testimonials02.text +\n      '\")'\n  );\n}\n", "ps1": "# Add two new bulleted list items after \"Building About Us - Part 02\":\n#   \"Building Testimonials - Part 01\"\n#   \"Building Testimonials - Part 02\"   (typed as \"...Part 0\" + \"2\")\n# Both new paragraphs must keep the same list/paragraph formatting\n# (ListParagraph style, numId=2 numbering, shading, spacing, outline\n# level, Segoe UI run fonts) as the existing \"Building About Us\" items.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n\n$anchorText = $anchor.Range.Text.TrimEnd([char]13, [char]7)\nif ($anchorText -ne \"Building About Us - Part 02\") {\n    throw \"Unexpected document structure: last paragraph is not 'Building About Us - Part 02' (found '$anchorText')\"\n}\n\n# Insert a new paragraph mark right after the anchor paragraph; Word\n# carries the anchor's paragraph formatting (style/numbering/shading/\n# spacing/outline level/run fonts) onto the new paragraph automatically.\n$anchor.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$p1.Range.Text = \"Building Testimonials - Part 01\"\n\n# Insert the second new bullet right after the first one, same formatting.\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)\n\n# The trailing \"2\" of \"Part 02\" was typed as its own run in the source\n# document, so write the base text first and then append \"2\" separately.\n$p2.Range.Text = \"Building Testimonials - Part 0\"\n$endRng = $p2.Range\n$endRng.Collapse(0)  # wdCollapseEnd\n$endRng.InsertAfter(\"2\")\n\n$finalText = $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text.TrimEnd([char]13, [char]7)\nif ($finalText -ne \"Building Testimonials - Part 02\") {\n    throw \"Failed to build 'Building Testimonials - Part 02' paragraph (found '$finalText')\"\n}\n"}
